$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) updates - force text storage so numeric-looking
# strings (e.g. "310.27", "19.50") are not coerced into Excel numbers
# (which would silently drop meaningful trailing zeros).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '23.881.18'
$ws.Range("D2").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.27'
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3878'
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3828'
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '51.22'
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.338'
$ws.Range("D10").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.08428'
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '23.81'
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.991'
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.998'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001314'
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.662.27'
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '93.94'
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06973'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '19.50'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.939'
$ws.Range("D21").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.62'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '23.883.37'
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.443'
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.914'
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.92'
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '153.20'
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.381'
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '137.09'
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.710'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.484'
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.826.39'
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.08150'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.9892'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02907'
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.637'
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.2673'
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '10.50'
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.09091'
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.7540'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '13.42'
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.423'
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '16.65'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6919'
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.438'
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.093'
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.001'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.08264'
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '133.69'
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.223'
$ws.Range("D51").Style = "Normal"

# Column E (Volume 1h) updates - plain text values (padded % strings),
# already safe from Excel numeric auto-conversion.
$ws.Range("E2").Value = '  -1.78%  '
$ws.Range("E3").Value = '  -1.02%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("E5").Value = '  -0.57%  '
$ws.Range("E6").Value = '  -0.08%  '
$ws.Range("E7").Value = '  -2.07%  '
$ws.Range("E8").Value = '  -2.59%  '
$ws.Range("E9").Value = '  -1.42%  '
$ws.Range("E10").Value = '  -3.37%  '
$ws.Range("E11").Value = '  -0.12%  '
$ws.Range("E12").Value = '  -1.58%  '
$ws.Range("E13").Value = '  -2.27%  '
$ws.Range("E14").Value = '  -4.33%  '
$ws.Range("E15").Value = '  -0.07%  '
$ws.Range("E16").Value = '  -1.34%  '
$ws.Range("E17").Value = '  -0.37%  '
$ws.Range("E18").Value = '  -1.72%  '
$ws.Range("E19").Value = '  -0.58%  '
$ws.Range("E20").Value = '  -4.68%  '
$ws.Range("E21").Value = '  -0.75%  '
$ws.Range("E22").Value = '  -0.08%  '
$ws.Range("E23").Value = '  -1.36%  '
$ws.Range("E24").Value = '  -1.80%  '
$ws.Range("E25").Value = '  -3.63%  '
$ws.Range("E26").Value = '  -5.74%  '
$ws.Range("E27").Value = '  -2.44%  '
$ws.Range("E28").Value = '  -2.25%  '
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("E30").Value = '  -3.61%  '
$ws.Range("E31").Value = '  -3.78%  '
$ws.Range("E32").Value = '  -1.44%  '
$ws.Range("E33").Value = '  -1.31%  '
$ws.Range("E34").Value = '  -1.40%  '
$ws.Range("E35").Value = '  -6.43%  '
$ws.Range("E36").Value = '  -5.91%  '
$ws.Range("E37").Value = '  -3.14%  '
$ws.Range("E38").Value = '  -2.97%  '
$ws.Range("E39").Value = '  -5.82%  '
$ws.Range("E40").Value = '  -2.11%  '
$ws.Range("E41").Value = '  -1.71%  '
$ws.Range("E42").Value = '  -2.47%  '
$ws.Range("E43").Value = '  -1.07%  '
$ws.Range("E44").Value = '  +0.72%  '
$ws.Range("E45").Value = '  -2.07%  '
$ws.Range("E46").Value = '  -3.48%  '
$ws.Range("E47").Value = '  -0.77%  '
$ws.Range("E48").Value = '  +0.04%  '
$ws.Range("E49").Value = '  -1.75%  '
$ws.Range("E50").Value = '  -2.01%  '
$ws.Range("E51").Value = '  -3.01%  '
